$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "LongI"
$ws.Range("B7").Value = 11
$ws.Range("C7").Value = 1.75
$ws.Range("D7").Value = 1.75

$ws.Range("B7").Select()
